$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.371.02'
$ws.Range("E2").Value = '  -1.18%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.378.85'
$ws.Range("E3").Value = '  +5.19%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.48'
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.642'
$ws.Range("E6").Value = '  -1.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.58'
$ws.Range("E7").Value = '  +14.22%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.481'
$ws.Range("E9").Value = '  +7.07%  '
$ws.Range("E10").Value = '  +0.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.92'
$ws.Range("E11").Value = '  -2.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '27.23'
$ws.Range("E12").Value = '  +2.82%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.728.24'
$ws.Range("E13").Value = '  +5.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.106'
$ws.Range("E14").Value = '  +0.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.12'
$ws.Range("E15").Value = '  +3.30%  '
$ws.Range("E16").Value = '  +2.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.862'
$ws.Range("E17").Value = '  +2.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.375.83'
$ws.Range("E18").Value = '  +4.69%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.380.15'
$ws.Range("E19").Value = '  -1.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0996'
$ws.Range("E20").Value = '  +1.72%  '
$ws.Range("E21").Value = '  +3.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.68'
$ws.Range("E22").Value = '  +1.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '250.58'
$ws.Range("E23").Value = '  +0.44%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.66'
$ws.Range("E25").Value = '  +4.06%  '
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.33'
$ws.Range("E27").Value = '  -3.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.07'
$ws.Range("E28").Value = '  +1.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.70'
$ws.Range("E29").Value = '  +3.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '173.99'
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("E31").Value = '  +6.45%  '
$ws.Range("E32").Value = '  -5.86%  '
$ws.Range("E33").Value = '  +0.42%  '
$ws.Range("E34").Value = '  +1.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0693'
$ws.Range("E35").Value = '  +0.88%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.09'
$ws.Range("E36").Value = '  +2.75%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.67'
$ws.Range("E37").Value = '  +4.07%  '
$ws.Range("E38").Value = '  +7.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.70'
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0256'
$ws.Range("E40").Value = '  +0.14%  '
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("E42").Value = '  +2.92%  '
$ws.Range("E43").Value = '  +7.22%  '
$ws.Range("E44").Value = '  +10.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '100.59'
$ws.Range("E45").Value = '  +1.87%  '
$ws.Range("E46").Value = '  +2.13%  '
$ws.Range("E47").Value = '  +2.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0957'
$ws.Range("E48").Value = '  +0.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.449.40'
$ws.Range("E49").Value = '  -0.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.601.42'
$ws.Range("E50").Value = '  +5.28%  '
$ws.Range("E51").Value = '  -1.99%  '
